$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("AC2").Value = 8.199999999999999

# Row 3
$ws.Range("Q3").Value = 1.86
$ws.Range("X3").Value = 1000
$ws.Range("Y3").Value = 1000
$ws.Range("Z3").Value = 980
$ws.Range("AC3").Value = 1000
$ws.Range("AD3").Value = 1000
$ws.Range("AF3").Value = 1000
$ws.Range("AG3").Value = 1000
$ws.Range("AJ3").Value = 980
$ws.Range("AK3").Value = 980
$ws.Range("AL3").Value = 980
$ws.Range("AN3").Value = 1000

# Row 5
$ws.Range("G5").Value = 1.66
$ws.Range("J5").Value = 4
$ws.Range("K5").Value = 4.9
$ws.Range("M5").Value = 1.06
$ws.Range("S5").Value = 2.86
$ws.Range("V5").Value = 1.13
$ws.Range("W5").Value = 2.52

# Row 6
$ws.Range("G6").Value = 3.4
$ws.Range("H6").Value = 2.32
$ws.Range("I6").Value = 2.6
$ws.Range("J6").Value = 3.4
$ws.Range("K6").Value = 3.75
$ws.Range("O6").Value = 1.32
$ws.Range("V6").Value = 1.62
$ws.Range("W6").Value = 1.42
$ws.Range("AA6").Value = 980
$ws.Range("AE6").Value = 980
$ws.Range("AI6").Value = 980
$ws.Range("AK6").Value = 980
$ws.Range("AL6").Value = 980
$ws.Range("AN6").Value = 42

# Row 7
$ws.Range("G7").Value = 1.5
$ws.Range("W7").Value = 3
$ws.Range("AG7").Value = 1000
$ws.Range("AL7").Value = 980
$ws.Range("AN7").Value = 8.800000000000001

# Row 8
$ws.Range("G8").Value = 1.82
$ws.Range("H8").Value = 5.7
$ws.Range("I8").Value = 9.199999999999999
$ws.Range("M8").Value = 1.1
$ws.Range("N8").Value = 2.58
$ws.Range("O8").Value = 1.5
$ws.Range("R8").Value = 1.19
$ws.Range("S8").Value = 4.5
$ws.Range("T8").Value = 2.36
$ws.Range("U8").Value = 1.6
$ws.Range("V8").Value = 1.14
$ws.Range("W8").Value = 2.22
$ws.Range("X8").Value = 1000
$ws.Range("Y8").Value = 1000
$ws.Range("AB8").Value = 1000
$ws.Range("AC8").Value = 1000
$ws.Range("AF8").Value = 1000
$ws.Range("AG8").Value = 1000
$ws.Range("AJ8").Value = 1000
